$d = $word.ActiveDocument

# --- Change 1: merge the two bold runs "State streaming stops and starts randomly" + ":" into a single run ---
$d.Content.Find.Execute("State streaming stops and starts randomly:", $false, $false, $false, $false, $false, $true, 1, $false, "State streaming stops and starts randomly:", 2) | Out-Null

# --- Change 2: replace the trailing empty paragraph with the new narrative paragraphs ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">Although, the fix above dramatically increased the length of time that the streaming worked, it did not allow it to work indefinitely. By disabling all other devices on the same SPI bus as the RFM95, the performance increased dramatically. This suggests that noise/signal integrity issues from other devices are causing the issues. To remedy this, I decreased the bus frequency from 6MB/s to 1.5MB/s. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>In future revisions of Strelka, putting the RFM95 on its own SPI bus would be preferrable.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">With the other devices on the SPI 2 bus disabled, the streaming was able to operate uninterrupted for 1 hour. This very strongly suggests that the other devices on the SPI bus are causing issues for the RFM95. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Disabling the BMX055 allowed to system to run for 20 </w:t>
      </w:r>
      <w:r>
        <w:t>minutes</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> before stalling. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Disabling the </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ASM330 allowed the system to run for about 50 mins. This cause of this may be because the ASM330’s I2C was not disabled using </w:t>
      </w:r>
      <w:r>
        <w:t>asm330lhhx_i2c_interface_set(&amp;asm330-&gt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>dev_ctx</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, ASM330LHHX_I2C_DISABLE);</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> which may have meant that every time CS went low, the ASM330 thought it was in I2C mode and started doing some random stuff. </w:t>
      </w:r>
    </w:p>
'@
$lastPara.Range.InsertXML($xml)
